$d = $word.ActiveDocument

# --- Edit 1: "Now go to the scripts directory..." paragraph ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("Now go to the scripts directory inside the mmitss-az repository:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "text 1 not found" }

$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D

$ins1 = $d.Range($r1.Start, $r1.Start)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Now go to the </w:t></w:r><w:r><w:t>' + $quoteOpen + '</w:t></w:r><w:r><w:t>build/</w:t></w:r><w:r><w:t>scripts</w:t></w:r><w:r><w:t>' + $quoteClose + '</w:t></w:r><w:r><w:t xml:space="preserve"> directory inside the mmitss-az repository:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins1.InsertXML($xml1)

# delete the old single run text (find it again, it was pushed after the newly inserted runs)
$r1b = $d.Content
$found1b = $r1b.Find.Execute("Now go to the scripts directory inside the mmitss-az repository:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1b) { throw "old text 1 not found for deletion" }
$r1b.Delete()
Write-Output "edit1 done"

# --- Edit 2: "cd mmitss-az/scripts" command line ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("cd mmitss-az/scripts", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "text 2 not found" }

$ins2 = $d.Range($r2.Start, $r2.Start)
$courier = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $courier + '<w:t>cd mmitss-az/</w:t></w:r><w:r>' + $courier + '<w:t>build/</w:t></w:r><w:r>' + $courier + '<w:t>scripts</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins2.InsertXML($xml2)

$r2b = $d.Content
$found2b = $r2b.Find.Execute("cd mmitss-az/scripts", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2b) { throw "old text 2 not found for deletion" }
$r2b.Delete()
Write-Output "edit2 done"
